$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.299.33"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "1.733.42"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "219.68"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "24.12"
$ws.Range("E8").Value = "  +11.05%  "
$ws.Range("D9").Value = "0.269"
$ws.Range("E9").Value = "  +4.70%  "
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.978.96"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").Value = "1.738.89"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("D16").Value = "67.71"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "28.289.11"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "242.14"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "0.0₃0755"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "7.93"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +2.22%  "
$ws.Range("D23").Value = "9.69"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "149.64"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("D27").Value = "16.67"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +3.18%  "
$ws.Range("E31").Value = "  +2.57%  "
$ws.Range("D32").Value = "3.44"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "1.504.58"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "2.41"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").Value = "70.50"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "5.66"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "1.882.64"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("D46").Value = "0.799"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  +9.21%  "
$ws.Range("D48").Value = "91.03"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +5.53%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.22"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  +0.48%  "
